$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("store")

function Set-StoreRow($RowIndex, $ItemName, $Price, $Timestamp, $Phone, $Picture) {
    $ws.Cells.Item($RowIndex, 1).Value = $ItemName
    $ws.Cells.Item($RowIndex, 2).Value = $Price
    $ws.Cells.Item($RowIndex, 3).Value = $Timestamp

    # "Phone" column holds digit-only text (e.g. "1234"); force it to stay
    # text instead of being auto-coerced to a number, then reset the style
    # back to Normal so no extra number-format style gets introduced.
    $ws.Cells.Item($RowIndex, 4).NumberFormat = "@"
    $ws.Cells.Item($RowIndex, 4).Value = $Phone
    $ws.Cells.Item($RowIndex, 4).Style = "Normal"

    $ws.Cells.Item($RowIndex, 5).Value = $Picture
}

Set-StoreRow 29 "Buzz Big Doll" 125 "12/08/2024 11:41:49" "1234" "pics/buzz_toy.png"
Set-StoreRow 30 "Barbie Doll"   150 "12/08/2024 11:44:22" "1234" "pics/barbie_doll.png"
Set-StoreRow 31 "Barbie Doll"   150 "12/08/2024 11:44:24" "1234" "pics/barbie_doll.png"
Set-StoreRow 32 "Barbie Doll"   150 "12/08/2024 11:44:25" "1234" "pics/barbie_doll.png"
Set-StoreRow 33 "Barbie Doll"   150 "12/08/2024 11:44:26" "1234" "pics/barbie_doll.png"
Set-StoreRow 34 "Barbie Doll"   150 "12/08/2024 11:44:28" "1234" "pics/barbie_doll.png"
